# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume update described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 <-> Row 14 swap: Litecoin and WrappedEther traded places in the ranking.
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'

# Price (D) / Volume(1h) (E) refresh for every affected row.
# A leading apostrophe forces plain-number-looking price strings to stay text,
# matching how these cells were already stored (t="inlineStr") before the edit.
$ws.Range("D2").Value = '30.510.36'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.872.10'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''247.45'
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = '''0.4733'
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("D8").Value = '''0.2891'
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("D9").Value = '''0.06462'
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("D10").Value = '''21.95'
$ws.Range("E10").Value = '  +2.60%  '
$ws.Range("D11").Value = '''0.07698'
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D12").Value = '''0.7378'
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").Value = '1.871.47'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '''95.83'
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("D15").Value = '''5.156'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").Value = '''274.29'
$ws.Range("E16").Value = '  -0.90%  '
$ws.Range("D17").Value = '30.568.37'
$ws.Range("D18").Value = '''13.21'
$ws.Range("E18").Value = '  -3.02%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '''0.000007466'
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").Value = '2.110.67'
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = '''5.203'
$ws.Range("E23").Value = '  -2.06%  '
$ws.Range("D24").Value = '''6.153'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").Value = '''165.26'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '''9.165'
$ws.Range("E26").Value = '  -1.69%  '
$ws.Range("D27").Value = '''18.63'
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("D28").Value = '''1.897'
$ws.Range("E28").Value = '  -4.81%  '
$ws.Range("D29").Value = '''0.09943'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  -2.80%  '
$ws.Range("D33").Value = '''4.067'
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("D34").Value = '''0.04753'
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("D35").Value = '''1.115'
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("D36").Value = '''0.6898'
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '''0.01851'
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").Value = '''2.755'
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("D40").Value = '''6.241'
$ws.Range("E40").Value = '  -4.14%  '
$ws.Range("D41").Value = '''72.90'
$ws.Range("E41").Value = '  +2.79%  '
$ws.Range("D42").Value = '''1.961'
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").Value = '''0.4142'
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").Value = '''0.8341'
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").Value = '''101.05'
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("D47").Value = '''9.326'
$ws.Range("E47").Value = '  -1.30%  '
$ws.Range("D48").Value = '''35.24'
$ws.Range("E48").Value = '  -0.35%  '
$ws.Range("D49").Value = '''6.945'
$ws.Range("E49").Value = '  -3.17%  '
$ws.Range("D50").Value = '''912.54'
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("D51").Value = '''0.05658'
$ws.Range("E51").Value = '  +0.94%  '
